$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "2016-03-02 05:59:30"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "2016-03-02 05:59:40"
